$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new text value in A2 (this will create a new shared string entry)
$ws.Range("A2").Value = "This is my second change"

# Column A width change (was 13.28515625 -> now 23.42578125) to fit the longer text
$ws.Columns.Item(1).ColumnWidth = 23.42578125

# Leave the active selection on A6 as in the target sheetView
$ws.Range("A6").Select()
